# The source data added a new weekly price observation for Mango at
# "Vega Central Mapocho de Santiago". In the canonical row order the new
# record sorts to row 522, pushing the former rows 522-577 down by one
# (to 523-578) while keeping their content untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 522; everything below (old 522..577) shifts down
# to 523..578, carrying its values/styles with it - exactly like a manual
# Excel "Insert Sheet Rows" on that row.
$ws.Rows.Item(522).Insert()

# Populate the newly inserted row 522 with the new observation.
$ws.Cells.Item(522, 1).Value  = 9
$ws.Cells.Item(522, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(522, 3).Value  = "Metropolitana"
$ws.Cells.Item(522, 4).Value  = 44946
$ws.Cells.Item(522, 5).Value  = 13
$ws.Cells.Item(522, 6).Value  = "Fruta"
$ws.Cells.Item(522, 7).Value  = 100108
$ws.Cells.Item(522, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(522, 9).Value  = 100108002
$ws.Cells.Item(522, 10).Value = "Mango"
$ws.Cells.Item(522, 11).Value = "Sin especificar"
$ws.Cells.Item(522, 12).Value = "Primera"
$ws.Cells.Item(522, 13).Value = 580
$ws.Cells.Item(522, 14).Value = 6000
$ws.Cells.Item(522, 15).Value = 7000
$ws.Cells.Item(522, 16).Value = 6483
$ws.Cells.Item(522, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(522, 18).Value = "Perú"
$ws.Cells.Item(522, 19).Value = 1621
$ws.Cells.Item(522, 20).Value = 4
